$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ParticipantsTab): update the Cypher queries in B2 and C2.
# Only change is 'Unknown' -> 'unknown' in the WHERE clause; all other text is unchanged.
$b2 = @'
Match (f)<--(g:genomic_info)
WHERE g.platform in ['unknown']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.gender,'') as `Gender`,
    coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`
LIMIT 100
'@

$statQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.platform in ['unknown']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH DISTINCT samp,s,p,f
RETURN
    count(distinct s) AS Studies,
    count(distinct p) AS Participants,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Files`
'@

$b3 = @'
Match (f)<--(g:genomic_info)
WHERE g.platform in ['unknown']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$b4 = @'
Match (f)<--(g:genomic_info)
WHERE g.platform in ['unknown']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH DISTINCT p,s,samp,f
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
   ORDER By f.file_name LIMIT 100
'@

$ws.Range("B2").Value = $b2
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $b3
$ws.Range("C3").Value = $statQuery
$ws.Range("B4").Value = $b4
$ws.Range("C4").Value = $statQuery

# Update the saved selection/active cell from D2 to B5.
$ws.Range("B5").Select()
